$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.277.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.797.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +6.03%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "117.27"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "341.33"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.553"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +5.57%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.19"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0871"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +7.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.17"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.63"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.229.54"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.802.04"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.889"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "52.071.64"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.24"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +11.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.37"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.96"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0989"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.48"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.41"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.84"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +11.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.96"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.42%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.21"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.20%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.04"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.42"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.72"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.97%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.05"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.00"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.94%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +29.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0373"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +13.39%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.33%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "127.03"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.25"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.107.80"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.36"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.23%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.920"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +22.58%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.98%  "
